# "Still obsessing over final report formatting and COLORS!"
#
# Applies color/line formatting tweaks to the "Modeling Flow" diagram on
# slide 8 of the capstone deck:
#   - "Model" headline text gets the signature blue (2DC3F0) fill
#   - four arrow connectors switch from solid blue to dotted green (59B600)
#   - the two "Cosine Sim." labels switch from blue to green text
#   - the two dotted blue "U-Turn" feedback arrows become solid, filled blue

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

$blue  = 15778605   # RGB(0x2D, 0xC3, 0xF0) packed as BGR: R + G*256 + B*65536
$green = 46681      # RGB(0x59, 0xB6, 0x00) packed as BGR: R + G*256 + B*65536

# --- 1. "Model" title text -> blue fill -------------------------------
$modelShape = $s.Shapes.Item("Content Placeholder 2")
$modelShape.TextFrame.TextRange.Font.Color.RGB = $blue

# --- 2. Four triangle-tipped connectors: blue solid -> green dotted ---
$dottedGreenConnectors = @(
    "Straight Arrow Connector 15",
    "Straight Arrow Connector 50",
    "Straight Arrow Connector 53",
    "Straight Arrow Connector 55"
)
foreach ($name in $dottedGreenConnectors) {
    $conn = $s.Shapes.Item($name)
    $conn.Line.ForeColor.RGB = $green
    $conn.Line.DashStyle = 2   # msoLineSysDot
}

# --- 3. The two "Cosine Sim." labels: blue text -> green text ---------
$cosineLabels = @("TextBox 20", "TextBox 92")
foreach ($name in $cosineLabels) {
    $label = $s.Shapes.Item($name)
    $label.TextFrame.TextRange.Font.Color.RGB = $green
}

# --- 4. The two U-Turn feedback arrows: outline-only dotted -> solid filled blue
$uTurnArrows = @("U-Turn Arrow 67", "U-Turn Arrow 68")
foreach ($name in $uTurnArrows) {
    $arrow = $s.Shapes.Item($name)
    $arrow.Fill.Solid()
    $arrow.Fill.ForeColor.RGB = $blue
    $arrow.Line.ForeColor.RGB = $blue
    $arrow.Line.Weight = 4.5       # 57150 EMU
    $arrow.Line.DashStyle = 1      # msoLineSolid
}
